$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to the Tuna price rows: each entry is Row number followed by
# a hashtable of Column letter -> new value (only changed cells listed).
$updates = @(
  @{ Row = 2; Cells = @{ "D"=44637; "M"=200; "N"=14000; "O"=15000; "P"=14500; "S"=806 } },
  @{ Row = 3; Cells = @{ "D"=44637; "M"=240; "N"=10000; "O"=11000; "P"=10500; "S"=583 } },
  @{ Row = 4; Cells = @{ "D"=44294; "L"="Especial"; "M"=200; "N"=14500; "O"=15000; "P"=14750; "S"=819 } },
  @{ Row = 5; Cells = @{ "D"=44294; "L"="Primera"; "M"=240; "N"=12500; "O"=13000; "P"=12750; "S"=708 } },
  @{ Row = 6; Cells = @{ "D"=44294; "L"="Segunda"; "M"=240; "N"=10500; "P"=10750; "S"=597 } },
  @{ Row = 7; Cells = @{ "D"=44631; "M"=240; "N"=15000; "O"=16000; "P"=15500; "S"=861 } },
  @{ Row = 8; Cells = @{ "D"=44631; "M"=248; "N"=12000; "O"=13000; "P"=12516; "S"=695 } },
  @{ Row = 9; Cells = @{ "D"=44631; "L"="Segunda"; "N"=9000; "O"=10000; "P"=9500; "S"=528 } },
  @{ Row = 10; Cells = @{ "D"=44607; "L"="Primera"; "M"=300; "N"=11000; "O"=12000; "P"=11500; "S"=639 } },
  @{ Row = 11; Cells = @{ "D"=44607; "L"="Segunda"; "M"=240; "N"=9000; "O"=10000; "P"=9500; "S"=528 } },
  @{ Row = 12; Cells = @{ "D"=44606; "L"="Primera"; "N"=11500; "O"=12000; "P"=11750; "S"=653 } },
  @{ Row = 13; Cells = @{ "D"=44606; "L"="Segunda"; "M"=240; "N"=9500; "O"=10000; "P"=9750; "S"=542 } },
  @{ Row = 14; Cells = @{ "D"=44610; "L"="Primera"; "N"=13000; "O"=14000; "P"=13500; "S"=750 } },
  @{ Row = 15; Cells = @{ "D"=44610; "L"="Segunda"; "M"=200; "N"=11000; "P"=11500; "S"=639 } },
  @{ Row = 16; Cells = @{ "D"=44603; "L"="Especial"; "N"=14500; "O"=15000; "P"=14750; "S"=819 } },
  @{ Row = 17; Cells = @{ "D"=44636; "N"=14000; "P"=14500; "S"=806 } },
  @{ Row = 18; Cells = @{ "D"=44636; "L"="Primera"; "M"=200; "N"=10000; "O"=11000; "P"=10500; "S"=583 } },
  @{ Row = 19; Cells = @{ "D"=44609; "M"=240; "N"=13000; "O"=14000; "P"=13500; "S"=750 } },
  @{ Row = 20; Cells = @{ "D"=44609; "L"="Segunda"; "M"=240; "N"=11000; "O"=12000; "P"=11500; "S"=639 } },
  @{ Row = 21; Cells = @{ "D"=44595; "N"=15500; "O"=16000; "P"=15750; "S"=875 } },
  @{ Row = 22; Cells = @{ "D"=44685; "L"="Especial"; "M"=200; "N"=19000; "O"=20000; "P"=19500; "S"=1083 } },
  @{ Row = 23; Cells = @{ "D"=44685; "L"="Primera"; "M"=160; "N"=15000; "O"=16000; "P"=15500; "S"=861 } },
  @{ Row = 24; Cells = @{ "D"=44634; "L"="Especial"; "M"=200; "N"=14000; "O"=15000; "P"=14500; "S"=806 } },
  @{ Row = 25; Cells = @{ "D"=44634; "L"="Primera"; "N"=10000; "O"=11000; "P"=10500; "S"=583 } },
  @{ Row = 26; Cells = @{ "D"=44687; "L"="Especial"; "M"=100; "N"=18000; "O"=19000; "P"=18500; "S"=1028 } },
  @{ Row = 27; Cells = @{ "D"=44687; "L"="Primera"; "M"=100; "N"=14000; "O"=15000; "P"=14500; "S"=806 } },
  @{ Row = 28; Cells = @{ "D"=44295; "L"="Especial"; "M"=200; "N"=14500; "O"=15000; "P"=14750; "S"=819 } },
  @{ Row = 29; Cells = @{ "D"=44295; "L"="Primera"; "M"=200; "N"=12500; "O"=13000; "P"=12750; "S"=708 } },
  @{ Row = 30; Cells = @{ "D"=44295; "L"="Segunda"; "M"=240; "N"=10500; "O"=11000; "P"=10750; "S"=597 } },
  @{ Row = 31; Cells = @{ "D"=44630; "L"="Especial"; "M"=300; "N"=15000; "O"=16000; "P"=15500; "S"=861 } },
  @{ Row = 32; Cells = @{ "D"=44630; "M"=300; "N"=12000; "O"=13000; "S"=694 } },
  @{ Row = 33; Cells = @{ "D"=44630; "M"=240; "N"=9000; "O"=10000; "P"=9500; "S"=528 } }
)

foreach ($update in $updates) {
  $r = $update.Row
  foreach ($col in $update.Cells.Keys) {
    $ws.Range("$col$r").Value = $update.Cells[$col]
  }
}

